$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 6780.96
$ws.Range("I28").Value = 833.9048
$ws.Range("J28").Value = 38003
$ws.Range("K28").Value = 833.9048
$ws.Range("L28").Value = 38003
$ws.Range("M28").Value = -348.9048
$ws.Range("N28").Value = -38973
$ws.Range("H29").Value = 30
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = $null
$ws.Range("H38").Value = 648.4231
$ws.Range("I38").Value = 181
$ws.Range("J38").Value = 1531.3334
$ws.Range("K38").Value = 543
$ws.Range("L38").Value = 4594.0002
$ws.Range("M38").Value = -171
$ws.Range("N38").Value = -5338.0002
$ws.Range("H40").Value = 950.2857
$ws.Range("I40").Value = 736.2
$ws.Range("J40").Value = 1069.2222
$ws.Range("K40").Value = 736.2
$ws.Range("L40").Value = 1069.2222
$ws.Range("M40").Value = -561.2
$ws.Range("N40").Value = -1419.2222
$ws.Range("H58").Value = 732.06665
$ws.Range("I58").Value = 705
$ws.Range("J58").Value = 1111
$ws.Range("K58").Value = 2115
$ws.Range("L58").Value = 3333
$ws.Range("M58").Value = -1965
$ws.Range("N58").Value = -3633
$ws.Range("H64").Value = 2899.875
$ws.Range("I64").Value = 2822.4443
$ws.Range("K64").Value = 2822.4443
$ws.Range("M64").Value = -2574.4443
$ws.Range("H67").Value = 2899.875
$ws.Range("I67").Value = 2822.4443
$ws.Range("K67").Value = 2822.4443
$ws.Range("M67").Value = -1964.4443
$ws.Range("H74").Value = 3846.9
$ws.Range("J74").Value = 3688
$ws.Range("L74").Value = 3688
$ws.Range("N74").Value = -5560
$ws.Range("H77").Value = 3846.9
$ws.Range("J77").Value = 3688
$ws.Range("L77").Value = 18440
$ws.Range("N77").Value = -27800
$ws.Range("H87").Value = 25466.467
$ws.Range("J87").Value = 25466.467
$ws.Range("L87").Value = 25466.467
$ws.Range("N87").Value = -27962.467
$ws.Range("H90").Value = 25466.467
$ws.Range("J90").Value = 25466.467
$ws.Range("L90").Value = 76399.401
$ws.Range("N90").Value = -88879.401

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4667.4634
$ws.Range("I63").Value = 3338.9
$ws.Range("J63").Value = 5096.032
$ws.Range("K63").Value = 3338.9
$ws.Range("L63").Value = 5096.032
$ws.Range("M63").Value = -2652.9
$ws.Range("N63").Value = -6468.032
$ws.Range("H66").Value = 4667.4634
$ws.Range("I66").Value = 3338.9
$ws.Range("J66").Value = 5096.032
$ws.Range("K66").Value = 16694.5
$ws.Range("L66").Value = 25480.16
$ws.Range("M66").Value = -13262.5
$ws.Range("N66").Value = -32344.16
$ws.Range("H80").Value = 20856.143
$ws.Range("J80").Value = 20856.143
$ws.Range("L80").Value = 20856.143
$ws.Range("N80").Value = -22852.143
$ws.Range("H83").Value = 20856.143
$ws.Range("J83").Value = 20856.143
$ws.Range("L83").Value = 62568.429
$ws.Range("N83").Value = -72552.429
$ws.Range("H88").Value = 3268.5715
$ws.Range("I88").Value = 2020
$ws.Range("J88").Value = 4933.3335
$ws.Range("K88").Value = 2020
$ws.Range("L88").Value = 4933.3335
$ws.Range("M88").Value = -1614
$ws.Range("N88").Value = -5745.3335
$ws.Range("H91").Value = 3268.5715
$ws.Range("I91").Value = 2020
$ws.Range("J91").Value = 4933.3335
$ws.Range("K91").Value = 2020
$ws.Range("L91").Value = 4933.3335
$ws.Range("M91").Value = -616
$ws.Range("N91").Value = -7741.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 17750.111
$ws.Range("I82").Value = 11252.333
$ws.Range("J82").Value = 20999
$ws.Range("K82").Value = 11252.333
$ws.Range("L82").Value = 20999
$ws.Range("M82").Value = -10869.333
$ws.Range("N82").Value = -21765
$ws.Range("H85").Value = 17750.111
$ws.Range("I85").Value = 11252.333
$ws.Range("J85").Value = 20999
$ws.Range("K85").Value = 11252.333
$ws.Range("L85").Value = 20999
$ws.Range("M85").Value = -9926.333000000001
$ws.Range("N85").Value = -23651
$ws.Range("H86").Value = 2064.8
$ws.Range("I86").Value = 1994.2222
$ws.Range("J86").Value = 2700
$ws.Range("K86").Value = 1994.2222
$ws.Range("L86").Value = 2700
$ws.Range("M86").Value = -871.2221999999999
$ws.Range("N86").Value = -4946
$ws.Range("H89").Value = 2064.8
$ws.Range("I89").Value = 1994.2222
$ws.Range("J89").Value = 2700
$ws.Range("K89").Value = 9971.110999999999
$ws.Range("L89").Value = 13500
$ws.Range("M89").Value = -4355.110999999999
$ws.Range("N89").Value = -24732

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10437.267
$ws.Range("I31").Value = 1237
$ws.Range("J31").Value = 13065.914
$ws.Range("K31").Value = 1237
$ws.Range("L31").Value = 13065.914
$ws.Range("M31").Value = -942
$ws.Range("N31").Value = -13655.914
$ws.Range("H34").Value = 10437.267
$ws.Range("I34").Value = 1237
$ws.Range("J34").Value = 13065.914
$ws.Range("K34").Value = 1237
$ws.Range("L34").Value = 13065.914
$ws.Range("M34").Value = -1035
$ws.Range("N34").Value = -13469.914
$ws.Range("H50").Value = 19248.75
$ws.Range("J50").Value = 19248.75
$ws.Range("L50").Value = 19248.75
$ws.Range("N50").Value = -20498.75
$ws.Range("H59").Value = 18999.166
$ws.Range("J59").Value = 18999.166
$ws.Range("L59").Value = 18999.166
$ws.Range("N59").Value = -21289.166
$ws.Range("H68").Value = 24599.5
$ws.Range("J68").Value = 24599.5
$ws.Range("L68").Value = 24599.5
$ws.Range("N68").Value = -26097.5
$ws.Range("H71").Value = 24599.5
$ws.Range("J71").Value = 24599.5
$ws.Range("L71").Value = 73798.5
$ws.Range("N71").Value = -81286.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 107143020
$ws.Range("I23").Value = 166.66667
$ws.Range("J23").Value = 136363790
$ws.Range("K23").Value = 500.00001
$ws.Range("L23").Value = 409091370
$ws.Range("M23").Value = -265.00001
$ws.Range("N23").Value = -409091840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 3999.9412
$ws.Range("J46").Value = 3999.9412
$ws.Range("L46").Value = 3999.9412
$ws.Range("N46").Value = -4311.9412
$ws.Range("H80").Value = 1882441
$ws.Range("I80").Value = 2253051.2
$ws.Range("J80").Value = 400000
$ws.Range("K80").Value = 2253051.2
$ws.Range("L80").Value = 400000
$ws.Range("M80").Value = -2252053.2
$ws.Range("N80").Value = -401996
$ws.Range("H83").Value = 1882441
$ws.Range("I83").Value = 2253051.2
$ws.Range("J83").Value = 400000
$ws.Range("K83").Value = 11265256
$ws.Range("L83").Value = 2000000
$ws.Range("M83").Value = -11260264
$ws.Range("N83").Value = -2009984
$ws.Range("H98").Value = 39999
$ws.Range("J98").Value = 39999
$ws.Range("L98").Value = 39999
$ws.Range("N98").Value = -45989

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 17617.166
$ws.Range("J22").Value = 20940.6
$ws.Range("L22").Value = 20940.6
$ws.Range("N22").Value = -21530.6
$ws.Range("H27").Value = 17617.166
$ws.Range("J27").Value = 20940.6
$ws.Range("L27").Value = 20940.6
$ws.Range("N27").Value = -21154.6
$ws.Range("H46").Value = 414.2857
$ws.Range("J46").Value = 380
$ws.Range("L46").Value = 380
$ws.Range("N46").Value = -756
$ws.Range("H55").Value = 503.6875
$ws.Range("I55").Value = 261.54544
$ws.Range("K55").Value = 261.54544
$ws.Range("M55").Value = -88.54543999999999
$ws.Range("H68").Value = 1873.2444
$ws.Range("I68").Value = 1760.9565
$ws.Range("J68").Value = 1990.6364
$ws.Range("K68").Value = 1760.9565
$ws.Range("L68").Value = 1990.6364
$ws.Range("M68").Value = -1011.9565
$ws.Range("N68").Value = -3488.6364
$ws.Range("H71").Value = 1873.2444
$ws.Range("I71").Value = 1760.9565
$ws.Range("J71").Value = 1990.6364
$ws.Range("K71").Value = 8804.782499999999
$ws.Range("L71").Value = 9953.182000000001
$ws.Range("M71").Value = -5060.782499999999
$ws.Range("N71").Value = -17441.182
$ws.Range("H74").Value = 24500
$ws.Range("J74").Value = 30000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -31996
$ws.Range("H77").Value = 24500
$ws.Range("J77").Value = 30000
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -99984

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2391.1143
$ws.Range("I122").Value = 1899.091
$ws.Range("J122").Value = 3223.7693
$ws.Range("K122").Value = 5697.272999999999
$ws.Range("L122").Value = 9671.3079
$ws.Range("M122").Value = -3247.272999999999
$ws.Range("N122").Value = -14571.3079
